$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.104551076889038
$ws.Range("B1").Value = 2.513625383377075
$ws.Range("C1").Value = 6.267829895019531
$ws.Range("D1").Value = 2.191050291061401
$ws.Range("E1").Value = 1.262075662612915
